$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.06692121177911758
$ws.Cells.Item(2, 2).Value = 0.9856381416320801
$ws.Cells.Item(2, 3).Value = 0.006765027996152639
$ws.Cells.Item(2, 4).Value = 0.9990689158439636
$ws.Cells.Item(3, 1).Value = 0.01193761173635721
$ws.Cells.Item(3, 2).Value = 0.9977198839187622
$ws.Cells.Item(3, 3).Value = 0.004301086533814669
$ws.Cells.Item(3, 4).Value = 0.9991405010223389
$ws.Cells.Item(4, 1).Value = 0.00746058439835906
$ws.Cells.Item(4, 2).Value = 0.9982528686523438
$ws.Cells.Item(4, 3).Value = 0.001642810530029237
$ws.Cells.Item(4, 4).Value = 0.9992837905883789
$ws.Cells.Item(5, 1).Value = 0.00375356781296432
$ws.Cells.Item(5, 2).Value = 0.998637855052948
$ws.Cells.Item(5, 3).Value = 0.001239861128851771
$ws.Cells.Item(5, 4).Value = 0.9997851252555847
$ws.Cells.Item(6, 1).Value = 0.002435752656310797
$ws.Cells.Item(6, 2).Value = 0.9993188977241516
$ws.Cells.Item(6, 3).Value = 0.0008243308984674513
$ws.Cells.Item(6, 4).Value = 0.9997134804725647
$ws.Cells.Item(7, 1).Value = 0.001477956306189299
$ws.Cells.Item(7, 2).Value = 0.9996446371078491
$ws.Cells.Item(7, 3).Value = 0.0008707991801202297
$ws.Cells.Item(7, 4).Value = 0.9997851252555847
$ws.Cells.Item(8, 1).Value = 0.001884446479380131
$ws.Cells.Item(8, 2).Value = 0.9995262026786804
$ws.Cells.Item(8, 3).Value = 0.0006819201516918838
$ws.Cells.Item(8, 4).Value = 0.9997851252555847
$ws.Cells.Item(9, 1).Value = 0.0008426898275502026
$ws.Cells.Item(9, 2).Value = 0.9998223185539246
$ws.Cells.Item(9, 3).Value = 0.0002715594018809497
$ws.Cells.Item(9, 4).Value = 0.99992835521698
$ws.Cells.Item(10, 1).Value = 0.0005546758184209466
$ws.Cells.Item(10, 2).Value = 0.9998815655708313
$ws.Cells.Item(10, 3).Value = [double]"6.735308124916628E-05"
$ws.Cells.Item(10, 4).Value = 0.99992835521698
$ws.Cells.Item(11, 1).Value = 0.001589049701578915
$ws.Cells.Item(11, 2).Value = 0.9994965791702271
$ws.Cells.Item(11, 3).Value = 0.0005091200582683086
$ws.Cells.Item(11, 4).Value = 0.9998567700386047
$ws.Cells.Item(12, 1).Value = 0.001085387193597853
$ws.Cells.Item(12, 2).Value = 0.9996742606163025
$ws.Cells.Item(12, 3).Value = 0.0004117018543183804
$ws.Cells.Item(12, 4).Value = 0.9998567700386047
$ws.Cells.Item(13, 1).Value = 0.000843483314383775
$ws.Cells.Item(13, 2).Value = 0.9996742606163025
$ws.Cells.Item(13, 3).Value = 0.0002512461505830288
$ws.Cells.Item(13, 4).Value = 0.99992835521698
$ws.Cells.Item(14, 1).Value = 0.0005470951437018812
$ws.Cells.Item(14, 2).Value = 0.9998223185539246
$ws.Cells.Item(14, 3).Value = 0.0001126685674535111
$ws.Cells.Item(14, 4).Value = 0.99992835521698
$ws.Cells.Item(15, 1).Value = 0.0005544557934626937
$ws.Cells.Item(15, 2).Value = 0.9998223185539246
$ws.Cells.Item(15, 3).Value = 0.0001423069334123284
$ws.Cells.Item(15, 4).Value = 0.99992835521698
$ws.Cells.Item(16, 1).Value = 0.001040547271259129
$ws.Cells.Item(16, 2).Value = 0.9997631311416626
$ws.Cells.Item(16, 3).Value = [double]"4.065388566232286E-05"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(17, 1).Value = 0.0003343665739521384
$ws.Cells.Item(17, 2).Value = 0.9998815655708313
$ws.Cells.Item(17, 3).Value = [double]"4.579811502480879E-05"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(18, 1).Value = 0.000751927902456373
$ws.Cells.Item(18, 2).Value = 0.9997631311416626
$ws.Cells.Item(18, 3).Value = [double]"4.693620212492533E-05"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 1).Value = 0.000488954596221447
$ws.Cells.Item(19, 2).Value = 0.9998815655708313
$ws.Cells.Item(19, 3).Value = [double]"7.639109389856458E-05"
$ws.Cells.Item(19, 4).Value = 0.99992835521698
$ws.Cells.Item(20, 1).Value = 0.001326893456280231
$ws.Cells.Item(20, 2).Value = 0.9998223185539246
$ws.Cells.Item(20, 3).Value = 0.0004605779831763357
$ws.Cells.Item(20, 4).Value = 0.9998567700386047
$ws.Cells.Item(21, 1).Value = 0.0006992131820879877
$ws.Cells.Item(21, 2).Value = 0.9997631311416626
$ws.Cells.Item(21, 3).Value = 0.0005651351530104876
$ws.Cells.Item(21, 4).Value = 0.9998567700386047
$ws.Cells.Item(22, 1).Value = 0.0002352878800593317
$ws.Cells.Item(22, 2).Value = 0.9999111890792847
$ws.Cells.Item(22, 3).Value = 0.0001679845590842888
$ws.Cells.Item(22, 4).Value = 0.99992835521698
$ws.Cells.Item(23, 1).Value = 0.000159932766109705
$ws.Cells.Item(23, 2).Value = 0.9999407529830933
$ws.Cells.Item(23, 3).Value = 0.0001876299502328038
$ws.Cells.Item(23, 4).Value = 0.99992835521698
$ws.Cells.Item(24, 1).Value = [double]"4.648940011975355E-05"
$ws.Cells.Item(24, 2).Value = 1
$ws.Cells.Item(24, 3).Value = 0.0002980238059535623
$ws.Cells.Item(24, 4).Value = 0.99992835521698
$ws.Cells.Item(25, 1).Value = 0.0001171449694084004
$ws.Cells.Item(25, 2).Value = 0.9999703764915466
$ws.Cells.Item(25, 3).Value = 0.001571018015965819
$ws.Cells.Item(25, 4).Value = 0.9998567700386047
$ws.Cells.Item(26, 1).Value = 0.0009362092823721468
$ws.Cells.Item(26, 2).Value = 0.9998519420623779
$ws.Cells.Item(26, 3).Value = [double]"6.185134679981275E-06"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.0003593254077713937
$ws.Cells.Item(27, 2).Value = 0.9998815655708313
$ws.Cells.Item(27, 3).Value = [double]"6.388336714735487E-06"
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = [double]"6.598445179406554E-05"
$ws.Cells.Item(28, 2).Value = 1
$ws.Cells.Item(28, 3).Value = [double]"3.372182254679501E-05"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 1).Value = [double]"4.929038914269768E-05"
$ws.Cells.Item(29, 2).Value = 0.9999703764915466
$ws.Cells.Item(29, 3).Value = [double]"8.050939868553542E-06"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = [double]"1.807944863685407E-05"
$ws.Cells.Item(30, 2).Value = 1
$ws.Cells.Item(30, 3).Value = [double]"6.886514893267304E-05"
$ws.Cells.Item(30, 4).Value = 0.99992835521698
$ws.Cells.Item(31, 1).Value = 0.0002705386723391712
$ws.Cells.Item(31, 2).Value = 0.9998815655708313
$ws.Cells.Item(31, 3).Value = [double]"2.721280338846555E-07"
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.000224950781557709
$ws.Cells.Item(32, 2).Value = 0.9999111890792847
$ws.Cells.Item(32, 3).Value = [double]"4.066272140335059E-06"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.001462705433368683
$ws.Cells.Item(33, 2).Value = 0.9998223185539246
$ws.Cells.Item(33, 3).Value = [double]"4.417074251250597E-06"
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = 0.0001107857897295617
$ws.Cells.Item(34, 2).Value = 0.9999407529830933
$ws.Cells.Item(34, 3).Value = [double]"1.094548133551143E-05"
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 1).Value = 0.000382153142709285
$ws.Cells.Item(35, 2).Value = 0.9998815655708313
$ws.Cells.Item(35, 3).Value = 0.0004731800290755928
$ws.Cells.Item(35, 4).Value = 0.9998567700386047
$ws.Cells.Item(36, 1).Value = [double]"6.422615115297958E-05"
$ws.Cells.Item(36, 2).Value = 0.9999407529830933
$ws.Cells.Item(36, 3).Value = 0.0009791573975235224
$ws.Cells.Item(36, 4).Value = 0.9998567700386047
$ws.Cells.Item(37, 1).Value = 0.001331298612058163
$ws.Cells.Item(37, 2).Value = 0.9997335076332092
$ws.Cells.Item(37, 3).Value = [double]"6.477295301010599E-06"
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = 0.0001470191782573238
$ws.Cells.Item(38, 2).Value = 0.9999407529830933
$ws.Cells.Item(38, 3).Value = [double]"6.450233286159346E-06"
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = [double]"9.806954039959237E-06"
$ws.Cells.Item(39, 2).Value = 1
$ws.Cells.Item(39, 3).Value = [double]"6.950530405447353E-06"
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 1).Value = 0.0003854723472613841
$ws.Cells.Item(40, 2).Value = 0.9999407529830933
$ws.Cells.Item(40, 3).Value = 0.0007687984616495669
$ws.Cells.Item(40, 4).Value = 0.9998567700386047
$ws.Cells.Item(41, 1).Value = 0.0002586292685009539
$ws.Cells.Item(41, 2).Value = 0.9998519420623779
$ws.Cells.Item(41, 3).Value = 0.0009060441516339779
$ws.Cells.Item(41, 4).Value = 0.9998567700386047
$ws.Cells.Item(42, 1).Value = 0.0003291761677246541
$ws.Cells.Item(42, 2).Value = 0.9998519420623779
$ws.Cells.Item(42, 3).Value = [double]"2.678345117601566E-05"
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 1).Value = [double]"2.449284875183366E-05"
$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(43, 3).Value = 0.0005817172932438552
$ws.Cells.Item(43, 4).Value = 0.9998567700386047
$ws.Cells.Item(44, 1).Value = 0.0008639546576887369
$ws.Cells.Item(44, 2).Value = 0.9998223185539246
$ws.Cells.Item(44, 3).Value = 0.000195757980691269
$ws.Cells.Item(44, 4).Value = 0.9998567700386047
$ws.Cells.Item(45, 1).Value = 0.0002635006094351411
$ws.Cells.Item(45, 2).Value = 0.9999407529830933
$ws.Cells.Item(45, 3).Value = 0.0001219434852828272
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 1).Value = [double]"9.173290891340002E-05"
$ws.Cells.Item(46, 2).Value = 0.9999703764915466
$ws.Cells.Item(46, 3).Value = 0.0005956810782663524
$ws.Cells.Item(46, 4).Value = 0.9998567700386047
$ws.Cells.Item(47, 1).Value = 0.0002630272065289319
$ws.Cells.Item(47, 2).Value = 0.9999407529830933
$ws.Cells.Item(47, 3).Value = 0.0006144459475763142
$ws.Cells.Item(47, 4).Value = 0.9998567700386047
$ws.Cells.Item(48, 1).Value = 0.0005530201597139239
$ws.Cells.Item(48, 2).Value = 0.9998519420623779
$ws.Cells.Item(48, 3).Value = 0.0003765238216146827
$ws.Cells.Item(48, 4).Value = 0.9998567700386047
$ws.Cells.Item(49, 1).Value = [double]"4.152861947659403E-05"
$ws.Cells.Item(49, 2).Value = 0.9999703764915466
$ws.Cells.Item(49, 3).Value = 0.0004278490960132331
$ws.Cells.Item(49, 4).Value = 0.9998567700386047
$ws.Cells.Item(50, 1).Value = 0.0008439061930403113
$ws.Cells.Item(50, 2).Value = 0.9998815655708313
$ws.Cells.Item(50, 3).Value = [double]"9.676550689619035E-06"
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 1).Value = [double]"1.925290598592255E-05"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = 0.0002361008519073948
$ws.Cells.Item(51, 4).Value = 0.9998567700386047
